$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a numeric-looking string into a cell while preserving
# its original (General) style and Text cell-type, i.e. no auto cast to Number.
function Set-TextValue($addr, $value) {
    $rng = $ws.Range($addr)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = $origStyle
}

Set-TextValue "D2" "243.45"
Set-TextValue "D3" "21.55"
Set-TextValue "D4" "5.309"
Set-TextValue "D5" "0.05653"
Set-TextValue "D6" "3.380"
Set-TextValue "D7" "0.8087"
Set-TextValue "D8" "0.9529"
Set-TextValue "D9" "0.1442"
Set-TextValue "D10" "0.07409"
Set-TextValue "D11" "0.03158"
Set-TextValue "D12" "0.03050"
Set-TextValue "D13" "0.09259"
Set-TextValue "D14" "3.576"
Set-TextValue "D15" "0.001644"
Set-TextValue "D16" "0.04707"
Set-TextValue "D17" "0.0005828"
$ws.Range("E17").Value = "16OneONEWorstin24h"
Set-TextValue "D18" "0.006356"
Set-TextValue "D19" "0.004994"
Set-TextValue "D21" "0.0001504"
Set-TextValue "D22" "0.0003104"
Set-TextValue "D23" "3.771"
Set-TextValue "D24" "6.383"
Set-TextValue "D25" "2.098"
Set-TextValue "D26" "0.3284"
Set-TextValue "D40" "0.03979"
Set-TextValue "D41" "0.006979"
$ws.Range("E41").Value = "40KickTokenKICK"
Set-TextValue "D42" "0.003508"
Set-TextValue "D43" "0.1035"
Set-TextValue "D44" "0.007487"
Set-TextValue "D45" "0.00005951"
Set-TextValue "D47" "0.0005508"
$ws.Range("E47").Value = "46ACDXExchangeACXT"
Set-TextValue "D48" "0.6833"
Set-TextValue "D49" "0.03260"
$ws.Range("E49").Value = "48BOLOBOLOBestin24h"
Set-TextValue "D50" "0.00002105"
